$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D (Price) and E (Volume 1h) for rows with unchanged Coin/Link
$ws.Range("D2").Value = "26.604.65"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "1.743.46"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("D4").Value = "`'0.9998"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "`'247.09"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "`'1.000"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "`'0.4928"
$ws.Range("E7").Value = "  +2.35%  "
$ws.Range("D8").Value = "`'0.2679"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "`'0.06300"
$ws.Range("E9").Value = "  +1.07%  "
$ws.Range("D10").Value = "1.743.01"
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("D11").Value = "`'0.07059"
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("D12").Value = "`'15.76"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "`'0.6160"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "`'4.593"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").Value = "`'78.27"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "`'1.000"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "26.616.99"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "`'0.000007321"
$ws.Range("E18").Value = "  +5.32%  "
$ws.Range("D19").Value = "`'1.000"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "`'11.58"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").Value = "1.965.60"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").Value = "`'4.588"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").Value = "`'8.748"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("D24").Value = "`'5.279"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("D25").Value = "`'139.57"
$ws.Range("E25").Value = "  +2.20%  "
$ws.Range("D26").Value = "`'15.49"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D28").Value = "`'1.768"
$ws.Range("E28").Value = "  -1.64%  "
$ws.Range("D29").Value = "`'107.75"
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("D30").Value = "`'4.056"
$ws.Range("E30").Value = "  +1.70%  "
$ws.Range("D31").Value = "`'0.08059"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").Value = "`'3.742"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "`'0.04638"
$ws.Range("E33").Value = "  +1.70%  "

# Rows 34-51: coin list shifted up by one (Frax dropped); update B,C,D,E
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "`'2.612"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "`'1.021"
$ws.Range("E35").Value = "  +3.00%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "`'0.6392"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "`'2.067"
$ws.Range("E37").Value = "  +3.40%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "`'0.9009"
$ws.Range("E38").Value = "  -4.55%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "`'2.428"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "`'1.004"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "`'0.01506"
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "`'101.99"
$ws.Range("E42").Value = "  -5.38%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "`'5.438"
$ws.Range("E43").Value = "  -3.87%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "`'0.3933"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "`'6.911"
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "`'0.1185"
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "`'0.05401"
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "`'7.863"
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "`'30.60"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "`'1.268"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "`'0.3436"
$ws.Range("E51").Value = "  +0.11%  "
